## rotational vector to euler angle
##
## The "Picture 6" floating picture (wp:docPr id="6") gets a rotation
## applied to its drawingml transform (a:xfrm/@rot). OOXML stores
## rotation in 60,000ths of a degree, so 2702005 == 45.0334... degrees.
## Word's Shape.Rotation property takes plain (clockwise) degrees and
## is written straight through to a:xfrm/@rot on save, so we hand it
## the exact fractional-degree equivalent of the target integer value.

$d = $word.ActiveDocument

$targetName = "Picture 6"
$rotationDegrees = 2702005 / 60000  # a:xfrm rot is in 60,000ths of a degree

$found = $false
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $shape = $d.Shapes.Item($i)
    if ($shape.Name -eq $targetName) {
        $shape.Rotation = $rotationDegrees
        $found = $true
        Write-Output ("Set rotation on '" + $shape.Name + "' to " + $shape.Rotation + " degrees")
    }
}

if (-not $found) {
    Write-Output "WARNING: target shape not found"
}
